$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabella 13")

# Add row 19: ID 18, DESCRIZIONE "Anagrafico di Unione Civile"
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Anagrafico di Unione Civile"

# Add row 20: ID 19, DESCRIZIONE "di Contratto di Convivenza"
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "di Contratto di Convivenza"

# Update selection to match target (D23)
$ws.Range("D23").Select()
